$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata": bump the "Last Updated" timestamp by one minute.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(2, 1).Value2 = "29 Oct 2025, 10:25 AM"

# ---------------------------------------------------------------------------
# Sheet "Top Gainers": rows 45-56 refreshed (leaderboard re-ranked, SAPPHIRE
# drops from the top of the window to the bottom with new figures).
# Columns: B=Stock, C=Latest, D=Weekly, E=Monthly
# ---------------------------------------------------------------------------
$gainers = $wb.Worksheets.Item("Top Gainers")

$gainersData = @(
    @(45, "SANDUMA",   4.593,  2.1405, 30.2813),
    @(46, "LLOYDSENT", 4.5646, 1.8339, 11.234),
    @(47, "STAR",      4.5025, 4.4319, 3.662),
    @(48, "RECLTD",    4.4992, 3.4756, 3.4062),
    @(49, "NBCC",      4.4511, 3.1605, 7.6018),
    @(50, "GPPL",      4.4154, 3.4073, 5.0497),
    @(51, "HUDCO",     4.3201, 3.8924, 5.3884),
    @(52, "SGMART",    4.2736, 8.258900000000001, 2.5381),
    @(53, "MRPL",      4.2642, 9.7103, 20.0542),
    @(54, "JKIL",      4.1372, 2.9463, 1.7584),
    @(55, "SAMBHV",    4.1349, 2.624,  5.167),
    @(56, "SAPPHIRE",  4.1265, 1.7633, -0.7999000000000001)
)

foreach ($row in $gainersData) {
    $r = $row[0]
    $gainers.Cells.Item($r, 2).Value2 = $row[1]
    $gainers.Cells.Item($r, 3).Value2 = $row[2]
    $gainers.Cells.Item($r, 4).Value2 = $row[3]
    $gainers.Cells.Item($r, 5).Value2 = $row[4]
}

# ---------------------------------------------------------------------------
# Sheet "Top Losers": a few isolated value tweaks plus rows 53-62 refreshed
# (PRIVISCL jumps from the bottom of the window to the top with new figures,
# CANHLIFE's Weekly/Monthly figures become unavailable this run).
# ---------------------------------------------------------------------------
$losers = $wb.Worksheets.Item("Top Losers")

# Isolated "Weekly" (column D) corrections, ticker/other columns unchanged.
$losers.Cells.Item(18, 4).Value2 = 5.978
$losers.Cells.Item(48, 4).Value2 = -2.9654

# Row 53: PRIVISCL, now with fully numeric D/E (was N/A before).
$losers.Cells.Item(53, 2).Value2 = "PRIVISCL"
$losers.Cells.Item(53, 3).Value2 = -2.6288
$losers.Cells.Item(53, 4).Value2 = -2.1048
$losers.Cells.Item(53, 5).Value2 = 19.7451

# Row 54: CANHLIFE, now with D/E unavailable (N/A text).
$losers.Cells.Item(54, 2).Value2 = "CANHLIFE"
$losers.Cells.Item(54, 3).Value2 = -2.6148
$losers.Cells.Item(54, 4).Value2 = "N/A"
$losers.Cells.Item(54, 5).Value2 = "N/A"

$losersData = @(
    @(55, "GKENERGY",   -2.6122, -9.807700000000001, 23.2758),
    @(56, "SGFIN",      -2.592,  -0.06270000000000001, 11.7235),
    @(57, "ARVINDFASN", -2.549,  -2.9892, -4.4223),
    @(58, "EDELWEISS",  -2.5422, -3.3745, 8.5305),
    @(59, "SAMHI",      -2.5284, 1.8231,  2.8516),
    @(60, "UJJIVANSFB", -2.5201, 0.3845,  12.6645),
    @(61, "AMBER",      -2.5098, -0.1082, 2.763),
    @(62, "GRPLTD",     -2.4898, -5.9894, -5.4586)
)

foreach ($row in $losersData) {
    $r = $row[0]
    $losers.Cells.Item($r, 2).Value2 = $row[1]
    $losers.Cells.Item($r, 3).Value2 = $row[2]
    $losers.Cells.Item($r, 4).Value2 = $row[3]
    $losers.Cells.Item($r, 5).Value2 = $row[4]
}

# ---------------------------------------------------------------------------
# Sheet "1 Month Performance": a couple of standalone % Change corrections
# plus rows 60-64 refreshed (GRMOVER drops from the top of the window to the
# bottom with a new figure).
# Columns: B=Stock, C=% Change
# ---------------------------------------------------------------------------
$perf = $wb.Worksheets.Item("1 Month Performance")

$perf.Cells.Item(53, 3).Value2 = 22.3984
$perf.Cells.Item(74, 3).Value2 = 18.8057

$perfData = @(
    @(60, "CEATLTD",    20.0239),
    @(61, "ATL",        19.9362),
    @(62, "SUBROS",     19.834),
    @(63, "HITECHGEAR", 19.8096),
    @(64, "GRMOVER",    19.7859)
)

foreach ($row in $perfData) {
    $r = $row[0]
    $perf.Cells.Item($r, 2).Value2 = $row[1]
    $perf.Cells.Item($r, 3).Value2 = $row[2]
}
